# Apply updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (e.g. thousands separated by
# "." and no true decimal semantics). Force text format first so Excel does not
# reinterpret values that happen to look numeric (e.g. "549.90" -> 549.9).
$ws.Range('D2').Value = '61.842.71'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').Value = '2.493.83'
$ws.Range('E3').Value = '  -3.97%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.90'
$ws.Range('E5').Value = '  -4.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.99'
$ws.Range('E6').Value = '  -5.89%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('D9').Value = '2.494.17'
$ws.Range('E9').Value = '  -3.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('E10').Value = '  -10.00%  '
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('E12').Value = '  -7.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').Value = '  -6.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.99'
$ws.Range('E14').Value = '  -7.19%  '
$ws.Range('D15').Value = '2.944.46'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').Value = '61.732.78'
$ws.Range('E16').Value = '  -2.70%  '
$ws.Range('E17').Value = '  -8.89%  '
$ws.Range('D18').Value = '2.490.54'
$ws.Range('E18').Value = '  -3.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.07'
$ws.Range('E19').Value = '  -7.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.00'
$ws.Range('E20').Value = '  -6.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  -7.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '320.64'
$ws.Range('E22').Value = '  -6.10%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -5.71%  '
$ws.Range('E25').Value = '  -3.93%  '
$ws.Range('D26').Value = '0.0₃0991'
$ws.Range('E26').Value = '  -8.02%  '
$ws.Range('D27').Value = '2.622.24'
$ws.Range('E27').Value = '  -3.43%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '533.55'
$ws.Range('E29').Value = '  -7.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.35'
$ws.Range('E30').Value = '  -8.69%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.47'
$ws.Range('E31').Value = '  -5.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.61'
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('E33').Value = '  -6.83%  '
$ws.Range('E34').Value = '  -8.11%  '
$ws.Range('E35').Value = '  -9.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.84'
$ws.Range('E36').Value = '  -10.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.86'
$ws.Range('E37').Value = '  -10.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('E39').Value = '  -5.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.44'
$ws.Range('E40').Value = '  -6.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '143.96'
$ws.Range('E41').Value = '  -6.47%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  -9.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.34'
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.27'
$ws.Range('E45').Value = '  -8.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '148.63'
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('E47').Value = '  -8.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.84'
$ws.Range('E48').Value = '  -10.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0534'
$ws.Range('E49').Value = '  -9.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.587'
$ws.Range('E50').Value = '  -6.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0943'
$ws.Range('E51').Value = '  -5.69%  '

Write-Host "Applied cryptos list update"
